# New weekly price record was added on top of the existing "Zanahoria"
# series for "Vega Modelo de Temuco": insert a brand-new row at row 111
# (pushing the previous rows 111..209 down to 112..210) and populate the
# new row with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 111 - shifts existing data (rows 111-209) down to 112-210.
$ws.Rows.Item(111).Insert()

$ws.Cells.Item(111, 1).Value2 = 10
$ws.Cells.Item(111, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(111, 3).Value2 = "La Araucanía"
$ws.Cells.Item(111, 4).Value2 = 44484
$ws.Cells.Item(111, 5).Value2 = 9
$ws.Cells.Item(111, 6).Value2 = 100114013
$ws.Cells.Item(111, 7).Value2 = "Zanahoria"
$ws.Cells.Item(111, 8).Value2 = "Sin especificar"
$ws.Cells.Item(111, 9).Value2 = "Primera"
$ws.Cells.Item(111, 10).Value2 = 20
$ws.Cells.Item(111, 11).Value2 = 9000
$ws.Cells.Item(111, 12).Value2 = 9000
$ws.Cells.Item(111, 13).Value2 = 9000
$ws.Cells.Item(111, 14).Value2 = "`$/saco 20 kilos"
$ws.Cells.Item(111, 15).Value2 = "Región del Maule"
$ws.Cells.Item(111, 16).Value2 = 450
$ws.Cells.Item(111, 17).Value2 = 20
$ws.Cells.Item(111, 18).Value2 = "Hortaliza"
